$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two data rows (JD_005 / Senior Devops Engineer / Demo,
# and JD_006 / Junior Devops Engineer / Testing) so the table shrinks
# from 7 rows to 5 rows (header + 4 job records).
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()

# Fill in the LinkedIn_Posted / Resume_received style status columns
# (F, G) for the remaining job rows that did not yet have them set.
$ws.Range("F3").Value = "Created"
$ws.Range("G3").Value = "Yes"

$ws.Range("F4").Value = "Created"
$ws.Range("G4").Value = "Yes"

$ws.Range("F5").Value = "Created"
$ws.Range("G5").Value = "Yes"

$ws.Range("G7").Select() | Out-Null
